$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like valid numbers to Excel (e.g. "1.795")
# must have their number format forced to Text ("@") before assignment,
# otherwise Excel auto-converts them to a floating point number instead of
# keeping the literal display string (matching the original inline string cells).
$textForceCells = @(
    "D5", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply each cell update from the diff
$ws.Range("D2").Value = "27.470.38"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.825.12"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "312.49"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.3614"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.07205"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").Value = "0.8592"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").Value = "20.58"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "1.824.92"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Value = "5.385"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "6.482"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "0.06931"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "80.30"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "0.000008876"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "15.33"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "27.606.88"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "5.119"
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("D23").Value = "10.92"
$ws.Range("E23").Value = "  +5.54%  "
$ws.Range("D24").Value = "2.060.89"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "155.07"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "18.68"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "5.126"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").Value = "114.14"
$ws.Range("E29").Value = "  -5.27%  "
$ws.Range("D30").Value = "1.795"
$ws.Range("D31").Value = "0.08851"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "0.7463"
$ws.Range("E32").Value = "  -2.46%  "
$ws.Range("D33").Value = "2.975"
$ws.Range("D34").Value = "4.532"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").Value = "1.118"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "1.085"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").Value = "0.05276"
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("D39").Value = "0.01914"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").Value = "2.782"
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("D41").Value = "0.5060"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "0.1641"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").Value = "6.429"
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").Value = "8.341"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "10.42"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").Value = "105.76"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.06448"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.4673"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "1.610"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("E51").Value = "  -1.33%  "
